$wb = $excel.ActiveWorkbook

# --- CTSViewPage: reorder existing rows, add two new test rows ---
$ws = $wb.Worksheets.Item("CTSViewPage")

# Capture the current (pre-edit) values of the four data rows (rows 2-5)
# before we start overwriting them, so the reorder is based on original data.
# NOTE: use .Text (not .Value) to read strings back reliably through this
# COM shim.
$r2a = $ws.Cells.Item(2,1).Text   # .../search/v?loc=0...NCI-2014-01507...  (Basic)
$r2c = $ws.Cells.Item(2,3).Text
$r3a = $ws.Cells.Item(3,1).Text   # .../search/v?loc=1...NCI-2015-01918...  (Advanced)
$r3c = $ws.Cells.Item(3,3).Text
$r4a = $ws.Cells.Item(4,1).Text   # .../search/v?t=C7711...NCI-2016-01041... (Basic)
$r4c = $ws.Cells.Item(4,3).Text
$r5a = $ws.Cells.Item(5,1).Text   # .../search/v?t=C9145...NCI-2011-02840... (Advanced)
$r5c = $ws.Cells.Item(5,3).Text

$bVal = $ws.Cells.Item(2,2).Text  # "CTS View Page" (same for every data row)

# New row order: old-3, old-5, old-2, old-4
$ws.Cells.Item(2,1).Value = $r3a
$ws.Cells.Item(2,2).Value = $bVal
$ws.Cells.Item(2,3).Value = $r3c

$ws.Cells.Item(3,1).Value = $r5a
$ws.Cells.Item(3,2).Value = $bVal
$ws.Cells.Item(3,3).Value = $r5c

$ws.Cells.Item(4,1).Value = $r2a
$ws.Cells.Item(4,2).Value = $bVal
$ws.Cells.Item(4,3).Value = $r2c

$ws.Cells.Item(5,1).Value = $r4a
$ws.Cells.Item(5,2).Value = $bVal
$ws.Cells.Item(5,3).Value = $r4c

# Replace the old direct-link test row with a real CTS view query, and
# add a second "Custom" search-type test case.
$ws.Cells.Item(6,1).Value = "/about-cancer/treatment/clinical-trials/search/v?id=NCT03200340&r=1"
$ws.Cells.Item(6,2).Value = $bVal
$ws.Cells.Item(6,3).Value = "Custom"

$ws.Cells.Item(7,1).Value = "/about-cancer/treatment/clinical-trials/search/v?id=NCI-2016-00402&r=1"
$ws.Cells.Item(7,2).Value = $bVal
$ws.Cells.Item(7,3).Value = "Custom"

# --- DynamicListingPage: move the selection (no data change) ---
$wsDyn = $wb.Worksheets.Item("DynamicListingPage")
$wsDyn.Range("A6").Select()

# --- Make CTSViewPage the active sheet/tab, with the new selection ---
$ws.Activate()
$ws.Range("A8").Select()

$wb.Save()
